# combine_merge / combine_stats fix:
# A handful of FCS ("Championship Subdivision") rows had been mismatched
# against the wrong FBS team's teamrankings row, leaking that other team's
# name (col E) and its PLpG3/PTpP3/OPLpG3/OPTpP3 stats (cols H:K) onto the
# wrong row, with "xxx" used as a placeholder in a couple of spots.
# Blank those mismatched cells out to a single space so they no longer show
# bogus data (mirrors upstream combine_stats/combine_merge behaviour).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 5, 8, 33, 49, 66, 85, 123, 164, 199)
$cols = @("E", "H", "I", "J", "K")

foreach ($r in $rows) {
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = " "
    }
}
